$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Department_ID" primary-column value for the DEPARTMENTS row was
# corrected to the upper-cased "DEPARTMENT_ID" (matching the casing
# convention used for the other Primary Column entries, e.g. EMPLOYEE_ID).
$ws.Range("D5").Value2 = "DEPARTMENT_ID"

# Re-apply the "Normal" style (with the 10pt font already used by the data
# rows) to the data range. This collapses the redundant duplicate cell
# format that only differed by an unused alignment flag, so every data
# cell consistently references the same format record.
$ws.Range("A2:E5").Style = "Normal"
$ws.Range("A2:E5").Font.Size = 10

# Column D needs to widen slightly to fit the new, slightly wider
# "DEPARTMENT_ID" text.
$ws.Columns("D").ColumnWidth = 12.983072916666666

# Move/save the active selection to D6, matching where the cursor ended up
# after making the edit.
$ws.Range("D6").Select()
